$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The CasesTab query (cell B2) had its final RETURN clause trimmed: the
# trailing "Cohort" column (coalesce(co.cohort_description, '') AS `Cohort`)
# was removed, and the preceding "Response to Treatment" line no longer
# ends with a comma.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (co:cohort)<-[*]-(c)
WHERE diag.disease_term IN ['Pulmonary Neoplasms']
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Cells.Item(2, 2).Value2 = $casesQuery

# Restore the view: scroll back to the top of the sheet, select B2, and use
# a normal (100%) zoom level instead of the previous zoomed-out / scrolled
# view.
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$window.Zoom = 100
$ws.Range("B2").Select()
